# Apply the authoring changes described in the commit:
#  - update the check_exists_in_zip rule text for the MicroplasticImages
#    zip file (trailing space added inside the parentheses)
#  - move the active selection to E19 (one row below the last data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("image_explorer_rules")

# Update the rule formula/text stored in E18
$ws.Range("E18").Value = "check_exists_in_zip(MicroplasticImages )"

# Update the saved cursor/selection position for the sheet
$ws.Range("E19").Select() | Out-Null
